$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
  @{rowNum=2; colF=446; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_s9are.png'; colM=90.14285714285714; colN=75.22857142857143; colO=82.68571428571428; colP=35; colQ=10; colR=10; colS=10; colT=10; colU=10; colV=10},
  @{rowNum=3; colF=447; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_jz3kd.png'; colM=72.79411764705883; colN=51.64705882352941; colO=62.22058823529412; colP=34; colQ=6; colR=6; colS=6; colT=6; colU=6; colV=6},
  @{rowNum=4; colF=448; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_as3da.png'; colM=84.53125; colN=63; colO=73.765625; colP=32; colQ=9; colR=9; colS=9; colT=9; colU=9; colV=9},
  @{rowNum=5; colF=449; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_cv9qj.png'; colM=60.34375; colN=35.34375; colO=47.84375; colP=32; colQ=3; colR=3; colS=3; colT=3; colU=3; colV=3},
  @{rowNum=6; colF=450; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_uy1n4.png'; colM=76.30555555555556; colN=55.33333333333334; colO=65.81944444444444; colP=36; colQ=7; colR=7; colS=7; colT=7; colU=7; colV=7},
  @{rowNum=7; colF=451; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_ncr40.png'; colM=75.66666666666667; colN=54.27272727272727; colO=64.96969696969697; colP=33; colQ=6; colR=6; colS=6; colT=6; colU=6; colV=6},
  @{rowNum=8; colF=452; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_hfz8w.png'; colM=55.46153846153846; colN=27.28205128205128; colO=41.37179487179487; colP=39; colQ=2; colR=2; colS=2; colT=2; colU=2; colV=2},
  @{rowNum=9; colF=453; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_fhm45.png'; colM=76.75; colN=57.71875; colO=67.234375; colP=32; colQ=7; colR=7; colS=7; colT=7; colU=7; colV=7},
  @{rowNum=10; colF=454; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_j5rpx.png'; colM=72.24242424242425; colN=50; colO=61.12121212121212; colP=33; colQ=5; colR=5; colS=5; colT=5; colU=5; colV=5},
  @{rowNum=11; colF=455; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_0j24m.png'; colM=63.6969696969697; colN=35.75757575757576; colO=49.72727272727273; colP=33; colQ=3; colR=3; colS=3; colT=3; colU=3; colV=3},
  @{rowNum=12; colF=456; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_bwo9g.png'; colM=64.81818181818181; colN=42.36363636363637; colO=53.59090909090909; colP=33; colQ=4; colR=4; colS=4; colT=4; colU=4; colV=4},
  @{rowNum=13; colF=457; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_i2k07.png'; colM=64.25925925925925; colN=40.92592592592592; colO=52.59259259259259; colP=27; colQ=3; colR=3; colS=3; colT=3; colU=3; colV=3},
  @{rowNum=14; colF=458; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_a8wvq.png'; colM=86.25925925925925; colN=66.25925925925925; colO=76.25925925925925; colP=27; colQ=10; colR=10; colS=10; colT=10; colU=10; colV=10},
  @{rowNum=15; colF=459; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_57os5.png'; colM=82.70588235294117; colN=65.73529411764706; colO=74.22058823529412; colP=34; colQ=9; colR=9; colS=9; colT=9; colU=9; colV=9},
  @{rowNum=16; colF=460; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_oz18d.png'; colM=78.93939393939394; colN=61.03030303030303; colO=69.98484848484848; colP=33; colQ=8; colR=8; colS=8; colT=8; colU=8; colV=8},
  @{rowNum=17; colF=461; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_84s7n.png'; colM=11.03125; colN=2.90625; colO=6.96875; colP=32; colQ=1; colR=1; colS=1; colT=1; colU=1; colV=1},
  @{rowNum=18; colF=462; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_h1yyu.png'; colM=64.8529411764706; colN=46.61764705882353; colO=55.73529411764706; colP=34; colQ=4; colR=4; colS=4; colT=4; colU=4; colV=4},
  @{rowNum=19; colF=463; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_8fpog.png'; colM=85.41666666666667; colN=72.30555555555556; colO=78.86111111111111; colP=36; colQ=10; colR=10; colS=10; colT=10; colU=10; colV=10},
  @{rowNum=20; colF=464; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_l9t30.png'; colM=67.2; colN=43.14285714285715; colO=55.17142857142858; colP=35; colQ=4; colR=4; colS=4; colT=4; colU=4; colV=4},
  @{rowNum=21; colF=465; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_01w8b.png'; colM=78.91891891891892; colN=61.21621621621622; colO=70.06756756756756; colP=37; colQ=8; colR=8; colS=8; colT=8; colU=8; colV=8},
  @{rowNum=22; colF=466; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_t1cr9.png'; colM=73.66666666666667; colN=53.51515151515152; colO=63.59090909090909; colP=33; colQ=6; colR=6; colS=6; colT=6; colU=6; colV=6},
  @{rowNum=23; colF=467; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_7pgd2.png'; colM=78.59375; colN=57.84375; colO=68.21875; colP=32; colQ=8; colR=7; colS=7; colT=7; colU=7; colV=7},
  @{rowNum=24; colF=468; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_xti0z.png'; colM=81.40625; colN=61.4375; colO=71.421875; colP=32; colQ=8; colR=8; colS=8; colT=8; colU=8; colV=8},
  @{rowNum=25; colF=469; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_ensho.png'; colM=72.7948717948718; colN=54.56410256410256; colO=63.67948717948718; colP=39; colQ=6; colR=6; colS=6; colT=6; colU=6; colV=6},
  @{rowNum=26; colF=470; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_qmand.png'; colM=86.11764705882354; colN=71.02941176470588; colO=78.57352941176471; colP=34; colQ=10; colR=10; colS=10; colT=10; colU=10; colV=10},
  @{rowNum=27; colF=471; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_faly8.png'; colM=33.41176470588236; colN=19.23529411764706; colO=26.32352941176471; colP=34; colQ=1; colR=1; colS=1; colT=1; colU=1; colV=1},
  @{rowNum=28; colF=472; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_uwv6y.png'; colM=78.88888888888889; colN=59.30555555555556; colO=69.09722222222223; colP=36; colQ=8; colR=8; colS=8; colT=8; colU=8; colV=8},
  @{rowNum=29; colF=473; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_c0me7.png'; colM=68.4; colN=45.62857142857143; colO=57.01428571428572; colP=35; colQ=4; colR=4; colS=4; colT=4; colU=4; colV=4},
  @{rowNum=30; colF=474; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_uspja.png'; colM=54.90909090909091; colN=29.12121212121212; colO=42.01515151515152; colP=33; colQ=2; colR=2; colS=2; colT=2; colU=2; colV=2},
  @{rowNum=31; colF=475; colH=$null; colI=$null; colJ='catch'; colK='f'; colL='stimuli/catch_28.jpg'; colM=$null; colN=$null; colO=$null; colP=$null; colQ=$null; colR=$null; colS=$null; colT=$null; colU=$null; colV=$null},
  @{rowNum=32; colF=476; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_05flq.png'; colM=47.10344827586207; colN=25.72413793103448; colO=36.41379310344828; colP=29; colQ=1; colR=1; colS=1; colT=1; colU=1; colV=1},
  @{rowNum=33; colF=477; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_ifebc.png'; colM=84; colN=65.88235294117646; colO=74.94117647058823; colP=34; colQ=10; colR=9; colS=9; colT=9; colU=9; colV=9},
  @{rowNum=34; colF=478; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_xesl0.png'; colM=69.28571428571429; colN=47.35714285714285; colO=58.32142857142857; colP=28; colQ=5; colR=5; colS=5; colT=5; colU=5; colV=5},
  @{rowNum=35; colF=479; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_411xa.png'; colM=51.03030303030303; colN=28.93939393939394; colO=39.98484848484848; colP=33; colQ=2; colR=2; colS=2; colT=2; colU=2; colV=2},
  @{rowNum=36; colF=480; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_z5osu.png'; colM=71.42857142857143; colN=47.34285714285714; colO=59.38571428571429; colP=35; colQ=5; colR=5; colS=5; colT=5; colU=5; colV=5},
  @{rowNum=37; colF=481; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_u9f9l.png'; colM=77.78571428571429; colN=57.25; colO=67.51785714285714; colP=28; colQ=7; colR=7; colS=7; colT=7; colU=7; colV=7},
  @{rowNum=38; colF=482; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_c79r7.png'; colM=56.26470588235294; colN=34.26470588235294; colO=45.26470588235294; colP=34; colQ=2; colR=2; colS=2; colT=2; colU=2; colV=2},
  @{rowNum=39; colF=483; colH='kitchens'; colI=$null; colJ='new'; colK='f'; colL='stimuli/img_ua9bs.png'; colM=82; colN=62.23333333333333; colO=72.11666666666667; colP=30; colQ=9; colR=9; colS=9; colT=9; colU=9; colV=9},
  @{rowNum=40; colF=484; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_5949k.png'; colM=60.8; colN=39.2; colO=50; colP=35; colQ=3; colR=3; colS=3; colT=3; colU=3; colV=3},
  @{rowNum=41; colF=485; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_qbkdt.png'; colM=69.45714285714286; colN=50.02857142857143; colO=59.74285714285715; colP=35; colQ=5; colR=5; colS=5; colT=5; colU=5; colV=5},
  @{rowNum=42; colF=486; colH='kitchens'; colI='target'; colJ='old'; colK='j'; colL='stimuli/img_k3abb.png'; colM=35.54054054054054; colN=16.54054054054054; colO=26.04054054054054; colP=37; colQ=1; colR=1; colS=1; colT=1; colU=1; colV=1}
)

foreach ($item in $rowData) {
    $rn = $item.rowNum
    $ws.Cells.Item($rn, 6).Value = $item.colF
    $ws.Cells.Item($rn, 8).Value = $item.colH
    $ws.Cells.Item($rn, 9).Value = $item.colI
    $ws.Cells.Item($rn, 10).Value = $item.colJ
    $ws.Cells.Item($rn, 11).Value = $item.colK
    $ws.Cells.Item($rn, 12).Value = $item.colL
    $ws.Cells.Item($rn, 13).Value = $item.colM
    $ws.Cells.Item($rn, 14).Value = $item.colN
    $ws.Cells.Item($rn, 15).Value = $item.colO
    $ws.Cells.Item($rn, 16).Value = $item.colP
    $ws.Cells.Item($rn, 17).Value = $item.colQ
    $ws.Cells.Item($rn, 18).Value = $item.colR
    $ws.Cells.Item($rn, 19).Value = $item.colS
    $ws.Cells.Item($rn, 20).Value = $item.colT
    $ws.Cells.Item($rn, 21).Value = $item.colU
    $ws.Cells.Item($rn, 22).Value = $item.colV
}

Write-Host "done"